$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference text already present in the sheet (rows 2-4) so the shared
# string table entries get reused instead of duplicated.
$pradosham = $ws.Range("B2").Value()
$pradoshamDesc = $ws.Range("C2").Value()
$pournami = $ws.Range("B4").Value()
$pournamiDesc = $ws.Range("C4").Value()

$pradoshamColor = $ws.Range("B2").Font.Color()
$pradoshamDescColor = $ws.Range("C2").Font.Color()
$pournamiColor = $ws.Range("B4").Font.Color()
$pournamiDescColor = $ws.Range("C4").Font.Color()

function Add-EventRow($row, $dateSerial, $eventText, $eventColor, $descText, $descColor) {
    $ws.Range("A$row").Value = $dateSerial

    $ws.Range("B$row").Value = $eventText
    $ws.Range("B$row").Font.Color = $eventColor

    $ws.Range("C$row").Value = $descText
    $ws.Range("C$row").Font.Color = $descColor
}

Add-EventRow 5  46054 $pournami  $pournamiColor  $pournamiDesc  $pournamiDescColor
Add-EventRow 6  46067 $pradosham $pradoshamColor $pradoshamDesc $pradoshamDescColor
Add-EventRow 7  46082 $pradosham $pradoshamColor $pradoshamDesc $pradoshamDescColor
Add-EventRow 8  46097 $pradosham $pradoshamColor $pradoshamDesc $pradoshamDescColor
Add-EventRow 9  46111 $pradosham $pradoshamColor $pradoshamDesc $pradoshamDescColor
Add-EventRow 10 46084 $pournami  $pournamiColor  $pournamiDesc  $pournamiDescColor

# Give the first new date cell the same "date, no wrap" number format used
# elsewhere in the workbook, then fan that exact style out to the rest of
# the new date cells via copy/paste-formats so a single new cellXfs entry
# is reused instead of one being minted per cell.
$ws.Range("A5").NumberFormat = "mm-dd-yy"
$ws.Range("A5").Copy()
$ws.Range("A6:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6:C9").Select()
